$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 150, pushing all existing data
# (old rows 150-260) down to rows 152-262.
$ws.Rows("150:151").Insert()

# Fill the two newly inserted rows with the new data records.

# Row 150 : new "Especial" grade record (Piña, Ecuador, 10-unit box)
$ws.Range("A150").Value = 4
$ws.Range("B150").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C150").Value = "Los Lagos"
$ws.Range("D150").Value = 44729
$ws.Range("E150").Value = 10
$ws.Range("F150").Value = "Fruta"
$ws.Range("G150").Value = 100108
$ws.Range("H150").Value = "Tropicales y subtropicales"
$ws.Range("I150").Value = 100108005
$ws.Range("J150").Value = "Piña"
$ws.Range("K150").Value = "Caramelo"
$ws.Range("L150").Value = "Especial"
$ws.Range("M150").Value = 200
$ws.Range("N150").Value = 18000
$ws.Range("O150").Value = 19000
$ws.Range("P150").Value = 18500
$ws.Range("Q150").Value = "$/caja 10 unidades"
$ws.Range("R150").Value = "Ecuador"
$ws.Range("S150").Value = 1850
$ws.Range("T150").Value = 10

# Row 151 : new "Tercera" grade record (Piña, Ecuador, 16-unit box)
$ws.Range("A151").Value = 4
$ws.Range("B151").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C151").Value = "Los Lagos"
$ws.Range("D151").Value = 44729
$ws.Range("E151").Value = 10
$ws.Range("F151").Value = "Fruta"
$ws.Range("G151").Value = 100108
$ws.Range("H151").Value = "Tropicales y subtropicales"
$ws.Range("I151").Value = 100108005
$ws.Range("J151").Value = "Piña"
$ws.Range("K151").Value = "Caramelo"
$ws.Range("L151").Value = "Tercera"
$ws.Range("M151").Value = 300
$ws.Range("N151").Value = 20000
$ws.Range("O151").Value = 21000
$ws.Range("P151").Value = 20500
$ws.Range("Q151").Value = "$/caja 16 unidades"
$ws.Range("R151").Value = "Ecuador"
$ws.Range("S151").Value = 1281
$ws.Range("T151").Value = 16
